$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Intro Fix Cross" body text (row 2, column B):
#    "fixation cross" -> "fixation DOT"
$ws.Range("B2").Value = "In this experiment you will need to keep your focus on the center of the fixation DOT shown below. "

# 2) Append three new rows (28-30) in column B with new onboarding copy,
#    styled with a purple Helvetica 10pt font (matching the new cellXfs/font
#    added to styles.xml).
$ws.Range("B28").Value = "Welcome to the BCBL!! Let's get going!!!"
$ws.Range("B29").Value = "In this experiment you will always need to keep your eyes on the center dot (shown below)"
$ws.Range("B30").Value = "It''s not always easy to fixate the cross because there will be many ''gratings'' on screen. As you''re about to see now…"

$r28 = $ws.Range("B28")
$r28.Font.Name = "Helvetica"
$r28.Font.Color = 15736992

$r29 = $ws.Range("B29")
$r29.Font.Name = "Helvetica"
$r29.Font.Color = 15736992

$r30 = $ws.Range("B30")
$r30.Font.Name = "Helvetica"
$r30.Font.Color = 15736992

# 3) Move the active selection to the newly-added last cell, as in the saved workbook.
$ws.Range("B30").Select() | Out-Null
